$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.751.37"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "3.854.03"
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "456.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.52%  "

$ws.Range("E6").Value = "  +13.79%  "

$ws.Range("E7").Value = "  +3.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.745"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.90%  "

$ws.Range("E10").Value = "  -2.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000318"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -7.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.72"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.24%  "

$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").Value = "4.462.44"
$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.88"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.25%  "

$ws.Range("D16").Value = "3.921.46"
$ws.Range("E16").Value = "  +2.83%  "

$ws.Range("E17").Value = "  -0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.13"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.94%  "

$ws.Range("D20").Value = "67.810.02"
$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "428.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +6.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.23"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +6.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +10.20%  "

$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +14.75%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "37.41"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.81"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.47"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "740.73"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.77%  "

$ws.Range("E32").Value = "  +11.27%  "

$ws.Range("E33").Value = "  -0.81%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.29"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +12.86%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.162"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +6.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.41"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.17%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.57"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.44%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0476"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.69%  "

$ws.Range("E40").Value = "  +15.20%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.66"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +20.00%  "

$ws.Range("D43").Value = "0.0₃0691"
$ws.Range("E43").Value = "  -8.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.35"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.18%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("E46").Value = "  +5.11%  "

$ws.Range("E47").Value = "  +4.38%  "

$ws.Range("E48").Value = "  +5.54%  "

$ws.Range("E49").Value = "  +6.01%  "

$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.75"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.32%  "

